# Auto-generated edit script: reproduces the results table being
# shifted up-and-left by one row/column (B2:E13 -> A1:D12), preserving
# per-cell formatting (borders/fonts) and values, then clearing the
# now-stale trailing row/column and fixing the view selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy each source cell's formatting onto its new (shifted) home.
#    Destination cells are written in row-major order (top-left to
#    bottom-right); because every source cell sits one row and one
#    column further down/right than its destination, it is never
#    overwritten before it has been read.
$ws.Range("B2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("E6").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("E10").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("B12").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Write the values into their new homes.
$ws.Range("A1").Value = 'Algorithm'
$ws.Range("B1").Value = 'Training Error (%)'
$ws.Range("C1").Value = 'Test Error (%)'
$ws.Range("D1").Value = 'Training Time (sec)'
$ws.Range("A2").Value = 'J48 Decision Tree'
$ws.Range("B2").Value = 5.13
$ws.Range("C2").Value = 31.31
$ws.Range("D2").Value = 2.94
$ws.Range("A3").Value = 'RF: 100 J48 trees'
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 25.16
$ws.Range("D3").Value = 10.71
$ws.Range("A4").Value = 'RF: 300 J48 Trees'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 24.99
$ws.Range("D4").Value = 32.49
$ws.Range("A5").Value = 'Logistic Regression'
$ws.Range("B5").Value = 28.71
$ws.Range("C5").Value = 40.450000000000003
$ws.Range("D5").Value = 41.11
$ws.Range("A6").Value = 'Naive Bayes'
$ws.Range("B6").Value = 33.93
$ws.Range("C6").Value = 45.73
$ws.Range("D6").Value = 0.17
$ws.Range("A7").Value = 'Adaboost: 20 J48 trees'
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 24.69
$ws.Range("D7").Value = 48.99
$ws.Range("A8").Value = 'Adabost: 100 J48 trees'
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 22.87
$ws.Range("D8").Value = 237.88
$ws.Range("A9").Value = 'Logitboost: 10 decision stumps'
$ws.Range("B9").Value = 29.15
$ws.Range("C9").Value = 48.16
$ws.Range("D9").Value = 9.5399999999999991
$ws.Range("A10").Value = 'Logitboost: 100 decision stumps'
$ws.Range("B10").Value = 21.8
$ws.Range("C10").Value = 38.86
$ws.Range("D10").Value = 94.26
$ws.Range("A11").Value = 'Logitboost: 100 stumps w/ pruning'
$ws.Range("B11").Value = 21.31
$ws.Range("C11").Value = 39.26
$ws.Range("D11").Value = 29.04
$ws.Range("A12").Value = 'Logitboost: 25 M5P trees'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 23.62
$ws.Range("D12").Value = 332.71

# 3) Remove the old bottom row (13) and right column (E), which are
#    no longer part of the table, clearing both contents and formats.
$ws.Range("B13:E13").Clear()
$ws.Range("E2:E12").Clear()

# 4) Update the sheet selection to match: column A selected.
$ws.Range("A1:A1048576").Select()
